$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Update column F ("想去人数" / number of people interested) on the
# "展览" (sheet1), "演出" (sheet2), and "全部类型" (sheet4) worksheets.
# "本地生活" (sheet3) has no data rows and needs no changes.
# -----------------------------------------------------------------------

# 展览
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    6  = 184
    7  = 3782
    8  = 185
    9  = 113
    10 = 93
    11 = 76
    12 = 83
    13 = 669
    15 = 907
    17 = 227
    21 = 84
    22 = 3283
    23 = 5624
    25 = 17
    26 = 82
    27 = 509
    29 = 3208
    31 = 2410
    33 = 512
    35 = 181
    36 = 247
    38 = 101
    39 = 495
    40 = 874
    42 = 27
    43 = 460
    45 = 535
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# 演出
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(2, 6).Value = 88

# 全部类型
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    6  = 184
    7  = 3782
    8  = 185
    9  = 113
    10 = 93
    11 = 88
    12 = 76
    13 = 83
    14 = 669
    16 = 907
    18 = 227
    22 = 84
    23 = 3283
    24 = 5624
    26 = 17
    27 = 82
    28 = 509
    30 = 3208
    32 = 2410
    34 = 512
    36 = 181
    37 = 247
    39 = 101
    40 = 495
    41 = 874
    43 = 27
    44 = 460
    46 = 535
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
